# Case_5_108 vm_pu.xlsx update: Bus 0 (slack) setpoint lowered from 1.05 pu to
# 1.02 pu, with the resulting per-bus voltage-magnitude results for the 24
# snapshots (rows 2-25, one row per time step) recomputed to the 380 kV case.
# Column H has no data in this sheet (gap between bus blocks) and must stay empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,12

# row 2 (time step 0)
$data[0,0] = 1.02
$data[0,1] = 1.042429616583995
$data[0,2] = 1.049457249154486
$data[0,3] = 1.050228828920262
$data[0,4] = 1.060249651280139
$data[0,5] = 1
$data[0,6] = $null
$data[0,7] = 1.039097680877022
$data[0,8] = 1.04750548468758
$data[0,9] = 1.052214254234562
$data[0,10] = 1.052983687805741
$data[0,11] = 1.062976943930868

# row 3 (time step 1)
$data[1,0] = 1.02
$data[1,1] = 1.043401473473572
$data[1,2] = 1.050203754475487
$data[1,3] = 1.051068603394769
$data[1,4] = 1.061115303175655
$data[1,5] = 1
$data[1,6] = $null
$data[1,7] = 1.039264375099796
$data[1,8] = 1.048123665262534
$data[1,9] = 1.052773179295484
$data[1,10] = 1.053635796209766
$data[1,11] = 1.063656854833152

# row 4 (time step 2)
$data[2,0] = 1.02
$data[2,1] = 1.044030973681405
$data[2,2] = 1.050687164070496
$data[2,3] = 1.051612867765385
$data[2,4] = 1.061676206222317
$data[2,5] = 1
$data[2,6] = $null
$data[2,7] = 1.039371074807582
$data[2,8] = 1.048523693946497
$data[2,9] = 1.053134545177883
$data[2,10] = 1.054057979907647
$data[2,11] = 1.064096926758423

# row 5 (time step 3)
$data[3,0] = 1.02
$data[3,1] = 1.044295768638591
$data[3,2] = 1.050890476723762
$data[3,3] = 1.051841884308696
$data[3,4] = 1.06191219226649
$data[3,5] = 1
$data[3,6] = $null
$data[3,7] = 1.039415652461198
$data[3,8] = 1.04869187080648
$data[3,9] = 1.053286391388146
$data[3,10] = 1.054235518857297
$data[3,11] = 1.064281961423094

# row 6 (time step 4)
$data[4,0] = 1.02
$data[4,1] = 1.04434023782243
$data[4,2] = 1.050924618886333
$data[4,3] = 1.05188034933701
$data[4,4] = 1.061951826012877
$data[4,5] = 1
$data[4,6] = $null
$data[4,7] = 1.039423120873681
$data[4,8] = 1.048720108691821
$data[4,9] = 1.05331188277931
$data[4,10] = 1.054265331497562
$data[4,11] = 1.064313031187159

# row 7 (time step 5)
$data[5,0] = 1.02
$data[5,1] = 1.04403451128444
$data[5,2] = 1.050689880403127
$data[5,3] = 1.051615927080967
$data[5,4] = 1.061679358764613
$data[5,5] = 1
$data[5,6] = $null
$data[5,7] = 1.039371671552987
$data[5,8] = 1.048525941116403
$data[5,9] = 1.053136574438745
$data[5,10] = 1.054060351985298
$data[5,11] = 1.064099399090579

# row 8 (time step 6)
$data[6,0] = 1.02
$data[6,1] = 1.042757926224945
$data[6,2] = 1.049709456460889
$data[6,3] = 1.05051245275892
$data[6,4] = 1.060542042359917
$data[6,5] = 1
$data[6,6] = $null
$data[6,7] = 1.039154256348442
$data[6,8] = 1.047714395988748
$data[6,9] = 1.052403206178902
$data[6,10] = 1.053204023404193
$data[6,11] = 1.063206696534087

# row 9 (time step 7)
$data[7,0] = 1.02
$data[7,1] = 1.04051339800035
$data[7,2] = 1.047984729623434
$data[7,3] = 1.048574749625077
$data[7,4] = 1.05854389892831
$data[7,5] = 1
$data[7,6] = $null
$data[7,7] = 1.038762264257143
$data[7,8] = 1.046284584716024
$data[7,9] = 1.051108701003355
$data[7,10] = 1.051696845312642
$data[7,11] = 1.061634647816359

# row 10 (time step 8)
$data[8,0] = 1.02
$data[8,1] = 1.039020454730541
$data[8,2] = 1.046836953737094
$data[8,3] = 1.047287579180301
$data[8,4] = 1.057215897810547
$data[8,5] = 1
$data[8,6] = $null
$data[8,7] = 1.038495000693024
$data[8,8] = 1.045331598191677
$data[8,9] = 1.050244274791141
$data[8,10] = 1.050693327463664
$data[8,11] = 1.06058736853711

# row 11 (time step 9)
$data[9,0] = 1.02
$data[9,1] = 1.038374815080541
$data[9,2] = 1.046340456978481
$data[9,3] = 1.046731336892779
$data[9,4] = 1.05664184901796
$data[9,5] = 1
$data[9,6] = $null
$data[9,7] = 1.038377871890806
$data[9,8] = 1.044919010028748
$data[9,9] = 1.049869645326844
$data[9,10] = 1.050259108238989
$data[9,11] = 1.060134079555237

# row 12 (time step 10)
$data[10,0] = 1.02
$data[10,1] = 1.038135119136259
$data[10,2] = 1.046156112237948
$data[10,3] = 1.046524892106737
$data[10,4] = 1.056428771170416
$data[10,5] = 1
$data[10,6] = $null
$data[10,7] = 1.038334154962514
$data[10,8] = 1.044765766652612
$data[10,9] = 1.049730443213382
$data[10,10] = 1.050097867695635
$data[10,11] = 1.059965737623938

# row 13 (time step 11)
$data[11,0] = 1.02
$data[11,1] = 1.0381865291315
$data[11,2] = 1.046195651335829
$data[11,3] = 1.046569167571929
$data[11,4] = 1.056474470305264
$data[11,5] = 1
$data[11,6] = $null
$data[11,7] = 1.038343541889631
$data[11,8] = 1.044798637397144
$data[11,9] = 1.049760304695792
$data[11,10] = 1.050132452157661
$data[11,11] = 1.06000184617191

# row 14 (time step 12)
$data[12,0] = 1.02
$data[12,1] = 1.038354999216643
$data[12,2] = 1.046325217420455
$data[12,3] = 1.046714268659357
$data[12,4] = 1.05662423287936
$data[12,5] = 1
$data[12,6] = $null
$data[12,7] = 1.038374262518198
$data[12,8] = 1.044906342668656
$data[12,9] = 1.049858139810505
$data[12,10] = 1.05024577906345
$data[12,11] = 1.060120163743929

# row 15 (time step 13)
$data[13,0] = 1.02
$data[13,1] = 1.038458815573265
$data[13,2] = 1.046405057511067
$data[13,3] = 1.046803692581793
$data[13,4] = 1.05671652639123
$data[13,5] = 1
$data[13,6] = $null
$data[13,7] = 1.038393162690621
$data[13,8] = 1.044972704828634
$data[13,9] = 1.049918412919758
$data[13,10] = 1.050315609879881
$data[13,11] = 1.060193067088632

# row 16 (time step 14)
$data[14,0] = 1.02
$data[14,1] = 1.039063320954634
$data[14,2] = 1.04686991520542
$data[14,3] = 1.047324518680401
$data[14,4] = 1.057254016410272
$data[14,5] = 1
$data[14,6] = $null
$data[14,7] = 1.038502744655688
$data[14,8] = 1.045358981688416
$data[14,9] = 1.050269130915621
$data[14,10] = 1.05072215180555
$data[14,11] = 1.06061745594117

# row 17 (time step 15)
$data[15,0] = 1.02
$data[15,1] = 1.039442730000336
$data[15,2] = 1.047161642466965
$data[15,3] = 1.047651517603971
$data[15,4] = 1.057591434345611
$data[15,5] = 1
$data[15,6] = $null
$data[15,7] = 1.038571107579598
$data[15,8] = 1.045601300093061
$data[15,9] = 1.050489040228158
$data[15,10] = 1.050977248784326
$data[15,11] = 1.060883715572233

# row 18 (time step 16)
$data[16,0] = 1.02
$data[16,1] = 1.039664111412381
$data[16,2] = 1.04733184997195
$data[16,3] = 1.047842357642732
$data[16,4] = 1.057788339401595
$data[16,5] = 1
$data[16,6] = $null
$data[16,7] = 1.038610847215302
$data[16,8] = 1.045742646163314
$data[16,9] = 1.050617277964093
$data[16,10] = 1.051126072430407
$data[16,11] = 1.061039038544481

# row 19 (time step 17)
$data[17,0] = 1.02
$data[17,1] = 1.039739610044052
$data[17,2] = 1.047389894417234
$data[17,3] = 1.047907447285958
$data[17,4] = 1.057855494984452
$data[17,5] = 1
$data[17,6] = $null
$data[17,7] = 1.038624374438099
$data[17,8] = 1.045790842480591
$data[17,9] = 1.050660998316722
$data[17,10] = 1.051176822483625
$data[17,11] = 1.061092002718773

# row 20 (time step 18)
$data[18,0] = 1.02
$data[18,1] = 1.039402014870103
$data[18,2] = 1.047130337924483
$data[18,3] = 1.047616422625922
$data[18,4] = 1.057555222769185
$data[18,5] = 1
$data[18,6] = $null
$data[18,7] = 1.038563786875735
$data[18,8] = 1.04557530100954
$data[18,9] = 1.050465449305339
$data[18,10] = 1.05094987619235
$data[18,11] = 1.060855146538546

# row 21 (time step 19)
$data[19,0] = 1.02
$data[19,1] = 1.038305385562102
$data[19,2] = 1.046287061320813
$data[19,3] = 1.046671535346283
$data[19,4] = 1.056580127396633
$data[19,5] = 1
$data[19,6] = $null
$data[19,7] = 1.03836522185351
$data[19,8] = 1.04487462585978
$data[19,9] = 1.049829331125199
$data[19,10] = 1.050212405785318
$data[19,11] = 1.060085321335817

# row 22 (time step 20)
$data[20,0] = 1.02
$data[20,1] = 1.037616605453966
$data[20,2] = 1.045757301810711
$data[20,3] = 1.046078422044168
$data[20,4] = 1.055967911568892
$data[20,5] = 1
$data[20,6] = $null
$data[20,7] = 1.038239160979643
$data[20,8] = 1.044434143328695
$data[20,9] = 1.049429100925166
$data[20,10] = 1.049749006052155
$data[20,11] = 1.059601474592385

# row 23 (time step 21)
$data[21,0] = 1.02
$data[21,1] = 1.037981672789205
$data[21,2] = 1.046038094964528
$data[21,3] = 1.046392749740601
$data[21,4] = 1.056292376227829
$data[21,5] = 1
$data[21,6] = $null
$data[21,7] = 1.038306103234731
$data[21,8] = 1.044667645445168
$data[21,9] = 1.04964129645536
$data[21,10] = 1.049994636363207
$data[21,11] = 1.059857954095896

# row 24 (time step 22)
$data[22,0] = 1.02
$data[22,1] = 1.039420412042296
$data[22,2] = 1.047144482951809
$data[22,3] = 1.047632280205094
$data[22,4] = 1.057571584928854
$data[22,5] = 1
$data[22,6] = $null
$data[22,7] = 1.038567095204805
$data[22,8] = 1.045587048858303
$data[22,9] = 1.050476109125789
$data[22,10] = 1.050962244596169
$data[22,11] = 1.060868055598916

# row 25 (time step 23)
$data[23,0] = 1.02
$data[23,1] = 1.041093066077208
$data[23,2] = 1.048430259418182
$data[23,3] = 1.04907488233404
$data[23,4] = 1.059059751938384
$data[23,5] = 1
$data[23,6] = $null
$data[23,7] = 1.038864652280607
$data[23,8] = 1.046654190871922
$data[23,9] = 1.051443617341201
$data[23,10] = 1.052086268635497
$data[23,11] = 1.062040932901625

$ws.Range("B2:M25").Value = $data

